$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

# 1) Flip status from Pendente to Resolvido for the rows the work instruction was updated on
$resolvedRows = @(288, 302, 303, 310)
foreach ($r in $resolvedRows) {
    $ws.Cells.Item($r, 10).Value = "Resolvido"
}

# 2) Append the new incident rows (week 24, 16/06/2025 - 20/06/2025)
$newRows = @(
    @{Row=314; A='ITI'; B='Erick da Silva'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=338021; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=315; A='ITI'; B='Erick da Silva'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337991; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=316; A='ITI'; B='Erick da Silva'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337816; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=317; A='ITI'; B='Erick da Silva'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337883; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=318; A='ITI'; B='Erick da Silva'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=338125; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=319; A='ITI'; B='Erick da Silva'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337838; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=320; A='ITI'; B='Felipe Nascimento'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337814; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=321; A='ITI'; B='Guilherme Worel'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=366994; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=322; A='ITI'; B='Guilherme Worel'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=338087; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=323; A='ITI'; B='Guilherme Worel'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337923; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=324; A='ITI'; B='Jorgenaldo Reis'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337938; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=325; A='ITI'; B='Lourival Moizés'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337799; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=326; A='ITI'; B='Michel Pessoa'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=338024; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=327; A='ITI'; B='Sostenes Simões'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=338018; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=328; A='ITI'; B='Sostenes Simões'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=338060; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=329; A='ITI'; B='Sostenes Simões'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337738; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'},
    @{Row=330; A='ITI'; B='Sostenes Simões'; C=2025; D=24; E='16/06/2025'; F='20/06/2025'; G=337872; H='06/2025'; I='16/06/2025'; J='Pendente'; K='Emerson Simette'}
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
}
